$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# The "Ana Van Gulick" author (author #2) previously had a redundant
# "url_name" attribute row (row 4). Remove that row.
# -----------------------------------------------------------------
$ws.Rows.Item(4).Delete()

# -----------------------------------------------------------------
# Add a new author, "William Shakespeare" (author #3), with richer
# metadata (name, first_name, last_name, email, orcid_id). Insert 5
# new rows right after the "Ana Van Gulick" row (now row 4), copying
# its formatting so the new rows inherit the same cell styles.
# -----------------------------------------------------------------
$ws.Rows.Item(4).Copy()
$ws.Range("A5:A9").EntireRow.Insert()

$ws.Range("A5").Value = "authors"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "name"
$ws.Range("D5").Value = "William Shakespeare"

$ws.Range("A6").Value = "authors"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "first_name"
$ws.Range("D6").Value = "William"

$ws.Range("A7").Value = "authors"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "last_name"
$ws.Range("D7").Value = "Shakespeare"

$ws.Range("A8").Value = "authors"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "email"
$ws.Range("D8").Value = "thebard@hotmail.com"

$ws.Range("A9").Value = "authors"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = "orcid_id"
$ws.Range("D9").Value = "000-000-12345"

# Turn the new author's email address into a live mailto: hyperlink.
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:thebard@hotmail.com")

# Update the saved selection to match the author's final cursor position.
$ws.Range("H18").Select() | Out-Null
